$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before existing row 3 (pushes old rows 3..16 down to 6..19)
$ws.Rows("3:5").Insert()

# New row 3
$ws.Cells.Item(3, 1).Value = "5a3c970775c03807006362c0"
$ws.Cells.Item(3, 2).Value = "EP1120171222052423980"
$ws.Cells.Item(3, 3).Value = "5a1bce0899c71b0900f9665a"
$ws.Cells.Item(3, 4).Value = "orgOwnerTest.20171117001.a"
$ws.Cells.Item(3, 5).Value = "5a1bd30999c71b0900f9665b"
$ws.Cells.Item(3, 6).Value = "orgAgentTest.20171127002.b"

# New row 4
$ws.Cells.Item(4, 1).Value = "5a3c954d75c03807006362a2"
$ws.Cells.Item(4, 2).Value = "EP0620171222051701196"
$ws.Cells.Item(4, 3).Value = "5a1bce0899c71b0900f9665a"
$ws.Cells.Item(4, 4).Value = "orgOwnerTest.20171117001.a"
$ws.Cells.Item(4, 5).Value = "5a1bd30999c71b0900f9665b"
$ws.Cells.Item(4, 6).Value = "orgAgentTest.20171127002.b"

# New row 5
$ws.Cells.Item(5, 1).Value = "5a39cc876556100800cbd47d"
$ws.Cells.Item(5, 2).Value = "EP0820171220023551609"
$ws.Cells.Item(5, 3).Value = "5a1bce0899c71b0900f9665a"
$ws.Cells.Item(5, 4).Value = "orgOwnerTest.20171117001.a"
$ws.Cells.Item(5, 5).Value = "5a1bd30999c71b0900f9665b"
$ws.Cells.Item(5, 6).Value = "orgAgentTest.20171127002.b"

# Update the row-count cell (C1): 14 -> 17
$ws.Cells.Item(1, 3).Value = 17

$wb.Save()
